# Auto-generated script applying scheduled market-data refresh to Halicarnassus_Profits sheets.
# For each (sheet, row) touched by the refresh, columns H/I/J/K/L (price data) and M/N
# (computed profit) are overwritten with the newly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 1945
$ws.Range("J16").Value = 1945
$ws.Range("L16").Value = 1945
$ws.Range("N16").Value = -2405
$ws.Range("H21").Value = 4999.5
$ws.Range("I21").Value = 4999.5
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 4999.5
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -4531.5
$ws.Range("H23").Value = 4999.5
$ws.Range("I23").Value = 4999.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4999.5
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -4765.5
$ws.Range("H34").Value = 4088.8
$ws.Range("I34").Value = 4088.8
$ws.Range("K34").Value = 4088.8
$ws.Range("M34").Value = -3885.8
$ws.Range("H36").Value = 4088.8
$ws.Range("I36").Value = 4088.8
$ws.Range("K36").Value = 4088.8
$ws.Range("M36").Value = -3373.8
$ws.Range("H88").Value = 1264.6666
$ws.Range("I88").Value = 1264.6666
$ws.Range("K88").Value = 1264.6666
$ws.Range("M88").Value = -858.6666
$ws.Range("H91").Value = 1264.6666
$ws.Range("I91").Value = 1264.6666
$ws.Range("K91").Value = 1264.6666
$ws.Range("M91").Value = 139.3334
$ws.Range("H96").Value = 188.1
$ws.Range("I96").Value = 226.71428
$ws.Range("J96").Value = 98
$ws.Range("K96").Value = 680.14284
$ws.Range("L96").Value = 294
$ws.Range("M96").Value = 692.85716
$ws.Range("N96").Value = -3040
$ws.Range("H103").Value = 5883.8335
$ws.Range("I103").Value = 6767.6665
$ws.Range("K103").Value = 20302.9995
$ws.Range("M103").Value = -19716.9995
$ws.Range("H112").Value = 1536
$ws.Range("I112").Value = 2866.6667
$ws.Range("J112").Value = 1250.8572
$ws.Range("K112").Value = 8600.000100000001
$ws.Range("L112").Value = 3752.5716
$ws.Range("M112").Value = -7492.000100000001
$ws.Range("N112").Value = -5968.571599999999
$ws.Range("H132").Value = 8681.566999999999
$ws.Range("I132").Value = 7523.161
$ws.Range("K132").Value = 22569.483
$ws.Range("M132").Value = -20039.483

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3103.6667
$ws.Range("I102").Value = 869.375
$ws.Range("K102").Value = 869.375
$ws.Range("M102").Value = 752.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2173.6667
$ws.Range("I99").Value = 3005
$ws.Range("K99").Value = 3005
$ws.Range("M99").Value = -1507

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2050
$ws.Range("I16").Value = 2050
$ws.Range("K16").Value = 2050
$ws.Range("M16").Value = -1763
$ws.Range("H113").Value = 2050
$ws.Range("I113").Value = 2050
$ws.Range("K113").Value = 2050
$ws.Range("M113").Value = 120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3907.6155
$ws.Range("I18").Value = 1449.75
$ws.Range("K18").Value = 4349.25
$ws.Range("M18").Value = -4180.25
$ws.Range("H26").Value = 190
$ws.Range("I26").Value = 190
$ws.Range("K26").Value = 570
$ws.Range("M26").Value = -282
$ws.Range("H39").Value = 11226.917
$ws.Range("J39").Value = 11226.917
$ws.Range("L39").Value = 33680.751
$ws.Range("N39").Value = -34268.751
$ws.Range("H130").Value = 1000
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020
$ws.Range("H132").Value = 2682.75
$ws.Range("I132").Value = 1492.6
$ws.Range("K132").Value = 13433.4
$ws.Range("M132").Value = -10903.4
$ws.Range("H138").Value = 6727.0713
$ws.Range("I138").Value = 1925.8
$ws.Range("K138").Value = 5777.4
$ws.Range("M138").Value = -637.3999999999996
$ws.Range("H139").Value = 4030.8572
$ws.Range("I139").Value = 4619.6665
$ws.Range("J139").Value = 498
$ws.Range("K139").Value = 13858.9995
$ws.Range("L139").Value = 1494
$ws.Range("M139").Value = -8718.999500000002
$ws.Range("N139").Value = -11774

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 698
$ws.Range("J80").Value = 897.3333
$ws.Range("L80").Value = 897.3333
$ws.Range("N80").Value = -2893.3333
$ws.Range("H83").Value = 698
$ws.Range("J83").Value = 897.3333
$ws.Range("L83").Value = 4486.6665
$ws.Range("N83").Value = -14470.6665
$ws.Range("H132").Value = 201281
$ws.Range("I132").Value = 334137
$ws.Range("J132").Value = 1997
$ws.Range("K132").Value = 1002411
$ws.Range("L132").Value = 5991
$ws.Range("M132").Value = -999881
$ws.Range("N132").Value = -11051

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1583
$ws.Range("I22").Value = 1462
$ws.Range("J22").Value = 1679.8
$ws.Range("K22").Value = 1462
$ws.Range("L22").Value = 1679.8
$ws.Range("M22").Value = -1167
$ws.Range("N22").Value = -2269.8
$ws.Range("H25").Value = 2502.5
$ws.Range("I25").Value = 2502.5
$ws.Range("K25").Value = 2502.5
$ws.Range("M25").Value = -2272.5
$ws.Range("H27").Value = 1583
$ws.Range("I27").Value = 1462
$ws.Range("J27").Value = 1679.8
$ws.Range("K27").Value = 1462
$ws.Range("L27").Value = 1679.8
$ws.Range("M27").Value = -1355
$ws.Range("N27").Value = -1893.8
$ws.Range("H68").Value = 10000
$ws.Range("J68").Value = 10000
$ws.Range("L68").Value = 10000
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 10000
$ws.Range("J71").Value = 10000
$ws.Range("L71").Value = 50000
$ws.Range("N71").Value = -57488
$ws.Range("H136").Value = 3197.2
$ws.Range("I136").Value = 2662
$ws.Range("K136").Value = 7986
$ws.Range("M136").Value = -5436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3466.5557
$ws.Range("I2").Value = 3862.5
$ws.Range("J2").Value = 299
$ws.Range("K2").Value = 3862.5
$ws.Range("L2").Value = 299
$ws.Range("M2").Value = -3750.5
$ws.Range("N2").Value = -523
$ws.Range("H4").Value = 12581.125
$ws.Range("J4").Value = 99.5
$ws.Range("L4").Value = 99.5
$ws.Range("N4").Value = -325.5
$ws.Range("H62").Value = 9945
$ws.Range("I62").Value = 7224.25
$ws.Range("K62").Value = 7224.25
$ws.Range("M62").Value = -6600.25
$ws.Range("H65").Value = 9945
$ws.Range("I65").Value = 7224.25
$ws.Range("K65").Value = 36121.25
$ws.Range("M65").Value = -33001.25
$ws.Range("H81").Value = 447.5
$ws.Range("I81").Value = 447.5
$ws.Range("K81").Value = 895
$ws.Range("M81").Value = 166
$ws.Range("H84").Value = 447.5
$ws.Range("I84").Value = 447.5
$ws.Range("K84").Value = 4475
$ws.Range("M84").Value = 829
